{"js": "const searchResults = context.document.body.search(\"02/12/22 - delayed meeting due to taking time off for holidays\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"2022/12/02 - delayed meeting due to taking time off for holidays\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"02/12/22 - delayed meeting due to taking time off for holidays\"\n$find.Replacement.Text = \"2022/12/02 - delayed meeting due to taking time off for holidays\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
